$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force affected cells to Text format so values are stored as strings (matching inlineStr source)
$ws.Range("D2:E26").NumberFormat = "@"
$ws.Range("D39:E47").NumberFormat = "@"
$ws.Range("D49:E51").NumberFormat = "@"

# Updated crypto price / 1h volume change data
$ws.Range("D2").Value = "307.06"
$ws.Range("E2").Value = "2.06%"
$ws.Range("D3").Value = "36.18"
$ws.Range("E3").Value = "3.46%"
$ws.Range("D4").Value = "5.095"
$ws.Range("E4").Value = "1.57%"
$ws.Range("D5").Value = "0.08124"
$ws.Range("E5").Value = "2.85%"
$ws.Range("D6").Value = "1.950"
$ws.Range("E6").Value = "0.92%"
$ws.Range("D7").Value = "7.765"
$ws.Range("E7").Value = "0.24%"
$ws.Range("D8").Value = "0.9312"
$ws.Range("E8").Value = "0.93%"
$ws.Range("D9").Value = "0.1418"
$ws.Range("E9").Value = "20.23%"
$ws.Range("D10").Value = "0.1925"
$ws.Range("E10").Value = "5.21%"
$ws.Range("D11").Value = "0.09257"
$ws.Range("E11").Value = "0.01%"
$ws.Range("D12").Value = "0.03523"
$ws.Range("E12").Value = "-0.30%"
$ws.Range("D13").Value = "0.09834"
$ws.Range("E13").Value = "-0.38%"
$ws.Range("D14").Value = "0.001423"
$ws.Range("E14").Value = "2.06%"
$ws.Range("D15").Value = "0.005868"
$ws.Range("E15").Value = "0.53%"
$ws.Range("D16").Value = "3.605"
$ws.Range("E16").Value = "2.86%"
$ws.Range("D17").Value = "4.190"
$ws.Range("E17").Value = "4.16%"
$ws.Range("D18").Value = "2.974"
$ws.Range("E18").Value = "0.34%"
$ws.Range("D19").Value = "0.3439"
$ws.Range("E19").Value = "-0.13%"
$ws.Range("D20").Value = "0.1349"
$ws.Range("E20").Value = "3.08%"
$ws.Range("D21").Value = "4.881"
$ws.Range("E21").Value = "-3.17%"
$ws.Range("D22").Value = "0.2408"
$ws.Range("E22").Value = "0.39%"
$ws.Range("D23").Value = "0.04511"
$ws.Range("E23").Value = "0.31%"
$ws.Range("D24").Value = "0.001218"
$ws.Range("E24").Value = "0.34%"
$ws.Range("D25").Value = "0.004875"
$ws.Range("E25").Value = "6.68%"
$ws.Range("E26").Value = "-0.80%"
$ws.Range("D39").Value = "0.02004"
$ws.Range("E39").Value = "5.88%"
$ws.Range("D40").Value = "0.04935"
$ws.Range("E40").Value = "5.18%"
$ws.Range("D41").Value = "0.01074"
$ws.Range("E41").Value = "12.35%"
$ws.Range("D42").Value = "0.007654"
$ws.Range("E42").Value = "1.07%"
$ws.Range("D43").Value = "0.1383"
$ws.Range("E43").Value = "4.53%"
$ws.Range("D44").Value = "0.002101"
$ws.Range("E44").Value = "-0.48%"
$ws.Range("D45").Value = "0.009991"
$ws.Range("E45").Value = "-10.37%"
$ws.Range("D46").Value = "0.00006444"
$ws.Range("E46").Value = "7.43%"
$ws.Range("D47").Value = "0.00000000750"
$ws.Range("E47").Value = "-0.01%"
$ws.Range("D49").Value = "0.001190"
$ws.Range("E49").Value = "-8.74%"
$ws.Range("D50").Value = "0.00002100"
$ws.Range("E50").Value = "-0.01%"
$ws.Range("D51").Value = "0.0002000"
$ws.Range("E51").Value = "-0.01%"

# Restore default (Normal) style so only text content differs from the original
$ws.Range("D2:E26").Style = "Normal"
$ws.Range("D39:E47").Style = "Normal"
$ws.Range("D49:E51").Style = "Normal"
